# Remove the second slide (the waffle plot quadrant breakdown slide) from
# the presentation, leaving only the first slide behind.
$p = $ppt.ActivePresentation

$s = $p.Slides.Item(2)
$s.Delete()
